# Automation Setup.pptx edits:
# 1. Remove the empty "Subtitle 2" placeholder shape from slide 1
#    (title slide no longer needs the subtitle placeholder).
# 2. Nudge "Picture 3" on slide 4 slightly to the left
#    (x offset 5203767 EMU -> 5157112 EMU; y unchanged).

$p = $ppt.ActivePresentation

# --- 1) Delete the "Subtitle 2" placeholder on slide 1 ---
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shape = $s1.Shapes.Item($i)
    if ($shape.Name -eq "Subtitle 2") {
        $shape.Delete()
    }
}

# --- 2) Move the picture on slide 4 ---
# Shape.Left/Top are expressed in points (1 pt = 12700 EMU). The target
# offset is 5157112 EMU; add a tiny half-EMU bias before dividing so the
# point->EMU round trip (which truncates) lands back on the exact value.
$EMU_PER_POINT = 12700
$targetEmuX = 5157112
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shape = $s4.Shapes.Item($i)
    if ($shape.Name -eq "Picture 3") {
        $shape.Left = ($targetEmuX + 0.5) / $EMU_PER_POINT
    }
}
